# BowlingAlly workbook update:
#  - Sheet1: reword a few task-tracker rows, add two more timing rows, renumber task list,
#            and insert a new "design" planning row.
#  - Add a new "Design " worksheet laying out the class design (Player / frame / Game)
#    with color-coded attribute/type/method legend.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Sheet1 edits
# ---------------------------------------------------------------------------

$timeFmt = $ws1.Range("D6").NumberFormat

# New note under the existing one at I3
$ws1.Range("I4").Value = "NB: every think must be worked before every commit "

# Two additional "time taken" rows
$ws1.Range("C7").Value = 0.22500000000000001
$ws1.Range("C7").NumberFormat = $timeFmt
$ws1.Range("F7").Value = 12

$ws1.Range("C8").Value = 0.23611111111111113
$ws1.Range("C8").NumberFormat = $timeFmt
$ws1.Range("D8").Value = 0.22500000000000001
$ws1.Range("D8").NumberFormat = $timeFmt
$ws1.Range("F8").Value = 16
$ws1.Range("I8").Value = "Make a nemurcal design "

# The old E8 "nbre" counter and the now-merged task rows are no longer used
$ws1.Range("E8").ClearContents()
$ws1.Range("E10").ClearContents()
$ws1.Range("I10").ClearContents()

# Renumbered task list (shifted down a row, with the two "testes" tasks merged)
$ws1.Range("E9").Value = 2
$ws1.Range("I9").Value = "Create packages"

$ws1.Range("E11").Value = 3
$ws1.Range("I11").Value = "create testes environment for every foo function "

$ws1.Range("E12").Value = 4
$ws1.Range("I12").Value = "Create function to calcule score of a given game "

$ws1.Range("E13").Value = 5
$ws1.Range("I13").Value = "create setRoll function "

$ws1.Range("E14").Value = 6
$ws1.Range("I14").Value = "create function that check score "

$ws1.Range("E15").Value = 7
$ws1.Range("I15").Value = "create function that print frames and final score and score after every frame"

$ws1.Range("E16").Value = 8
$ws1.Range("I16").Value = "add multiple player feature "

$ws1.Range("E17").Value = 9
$ws1.Range("I17").Value = "add rank after gameover"

$ws1.Range("E18").Value = 10
$ws1.Range("I18").Value = "calcule score after every roll setted in the game "

$ws1.Range("E19").Value = 11
$ws1.Range("I19").Value = "next feature pins presentation feature "

$ws1.Range("E20").Value = 12

$ws1.Range("E11:E20").Select()

# ---------------------------------------------------------------------------
# Add the "Design " worksheet (after Sheet1) and make it the active sheet
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Design "

# Font colors (BGR-encoded OLE colors, matching the workbook's new palette)
$blue = 15773696    # RGB FF00B0F0 - class-name banner
$green = 5287936    # RGB FF00B050 - "Types" legend / *_type names
$purple = 10498160  # RGB FF7030A0 - attribute / method values
$accent1 = 5         # xlThemeColorAccent1 - theme-based blue used once on the "Game" type cell

$ws2.Range("E2").Value = "Classes"

# --- Player class -----------------------------------------------------
$ws2.Range("F3").Value = "class name "
$ws2.Range("G3").Value = "Player "
$ws2.Range("G3").Font.Color = $blue

$ws2.Range("F4").Value = "atribute"
$ws2.Range("G4").Value = "type"
$ws2.Range("H4").Value = "methodes"
$ws2.Range("I4").Value = "type"

$ws2.Range("F5").Value = "Name "
$ws2.Range("F5").Font.Color = $purple
$ws2.Range("G5").Value = "String "
$ws2.Range("H5").Value = "None"
$ws2.Range("I5").Value = "None"

$ws2.Range("F6").Value = "Game "
$ws2.Range("F6").Font.Color = $purple
$ws2.Range("G6").Value = "Game"
$ws2.Range("G6").Font.ThemeColor = $accent1

$ws2.Range("F7").Value = "score "
$ws2.Range("F7").Font.Color = $purple
$ws2.Range("G7").Value = "Integer "

$ws2.Range("F8").Value = "Current round"
$ws2.Range("F8").Font.Color = $purple
$ws2.Range("G8").Value = "Integer "

# --- frame class --------------------------------------------------------
$ws2.Range("J3").Value = "class name "
$ws2.Range("K3").Value = "frame  "
$ws2.Range("K3").Font.Color = $blue

$ws2.Range("J4").Value = "atribute"
$ws2.Range("K4").Value = "type"
$ws2.Range("L4").Value = "methodes"
$ws2.Range("M4").Value = "type"

$ws2.Range("J5").Value = "Roll1"
$ws2.Range("J5").Font.Color = $purple
$ws2.Range("K5").Value = "Roll_type"
$ws2.Range("K5").Font.Color = $green
$ws2.Range("L5").Value = "SetRoll"
$ws2.Range("L5").Font.Color = $purple

$ws2.Range("J6").Value = "Roll2"
$ws2.Range("J6").Font.Color = $purple
$ws2.Range("K6").Value = "Roll_type"
$ws2.Range("K6").Font.Color = $green

# --- Game class ----------------------------------------------------------
$ws2.Range("N3").Value = "class name "
$ws2.Range("O3").Value = "Game"
$ws2.Range("O3").Font.Color = $blue

$ws2.Range("N4").Value = "atribute"
$ws2.Range("O4").Value = "type"
$ws2.Range("P4").Value = "methodes"

$ws2.Range("N5").Value = "Game"
$ws2.Range("N5").Font.Color = $purple
$ws2.Range("O5").Value = "Game_type"
$ws2.Range("O5").Font.Color = $green
$ws2.Range("P5").Value = "Get previous game Type "
$ws2.Range("P5").Font.Color = $purple

$ws2.Range("P6").Value = "Get actual frame number "
$ws2.Range("P6").Font.Color = $purple

$ws2.Range("P7").Value = "get Actual frame number "
$ws2.Range("P7").Font.Color = $purple

# --- Objects / Types legend (left column) --------------------------------
$ws2.Range("A4").Value = "Objects "
$ws2.Range("A4").Font.Color = $blue

$ws2.Range("A5").Value = "Types"
$ws2.Range("A5").Font.Color = $green

$ws2.Range("A6").Value = "attribute"
$ws2.Range("A6").Font.Color = $purple

$ws2.Range("A7").Value = "methodes"
$ws2.Range("A7").Font.Color = $purple

# --- New types legend ------------------------------------------------------
$ws2.Range("E13").Value = "New types"

$ws2.Range("F14").Value = "Game_type"
$ws2.Range("F14").Font.Color = $green
$ws2.Range("G14").Value = "array of 10 frames "

$ws2.Range("F15").Value = "Roll_type"
$ws2.Range("F15").Font.Color = $green
$ws2.Range("G15").Value = "Integer between 0 and 10 "

# --- column widths (approximate best-fit sizing) --------------------------
$offset = 5.0 / 6.0
$ws2.Columns.Item(6).ColumnWidth = 12.33203125 - $offset
$ws2.Columns.Item(7).ColumnWidth = 10.77734375 - $offset
$ws2.Columns.Item(9).ColumnWidth = 5.44140625 - $offset
$ws2.Columns.Item(10).ColumnWidth = 10.33203125 - $offset
$ws2.Columns.Item(11).ColumnWidth = 8.6640625 - $offset
$ws2.Columns.Item(12).ColumnWidth = 9.109375 - $offset
$ws2.Columns.Item(13).ColumnWidth = 4.5546875 - $offset
$ws2.Columns.Item(14).ColumnWidth = 10.33203125 - $offset
$ws2.Columns.Item(15).ColumnWidth = 16.5546875 - $offset
$ws2.Columns.Item(16).ColumnWidth = 22.109375 - $offset
$ws2.Columns.Item(17).ColumnWidth = 16.5546875 - $offset
$ws2.Columns.Item(18).ColumnWidth = 22.109375 - $offset

$ws2.Range("T22").Select()

Write-Host "done"
